# This script updates the cryptocurrency price/volume table (and a few
# coin name / link cells for rows that shifted position in the source
# feed) on the active worksheet, matching the "Updated cryptos list"
# GitHub Actions commit.
#
# Columns B and C hold plain text (coin name / coinranking.com link).
# Columns D and E hold values that LOOK numeric (e.g. "1.000", "240.30",
# "  -0.35%  ") but must stay stored as literal text, exactly as they
# were in the original workbook (t="inlineStr" / General-format text
# cells). Assigning such strings straight to Range.Value lets Excel's
# COM layer auto-parse/convert them into real numbers, which would
# corrupt values like "1.000" -> 1. To avoid that we force the cell to
# Text format before writing the value and then restore the "Normal"
# cell style afterwards (so no stray formatting is left behind), which
# keeps Excel from reinterpreting the string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = [ordered]@{
    'D2' = '29.360.14'
    'E2' = '  -0.35%  '
    'D3' = '1.847.62'
    'E3' = '  -0.26%  '
    'D4' = '1.000'
    'E4' = '  -0.07%  '
    'D5' = '240.30'
    'E5' = '  -0.20%  '
    'D6' = '0.6281'
    'E6' = '  -0.50%  '
    'D7' = '1.001'
    'E7' = '  -0.08%  '
    'D8' = '0.07565'
    'E8' = '  -1.36%  '
    'D9' = '0.2906'
    'E9' = '  -1.10%  '
    'D10' = '24.55'
    'E10' = '  -0.51%  '
    'D11' = '0.07749'
    'D12' = '1.847.06'
    'E12' = '  -0.68%  '
    'D13' = '5.012'
    'E13' = '  -0.50%  '
    'D14' = '0.6783'
    'E14' = '  -0.38%  '
    'D15' = '0.00001039'
    'E15' = '  -2.05%  '
    'D16' = '83.06'
    'E16' = '  -0.71%  '
    'D17' = '6.107'
    'E17' = '  -1.06%  '
    'D18' = '29.379.38'
    'E18' = '  -0.38%  '
    'D19' = '228.88'
    'E19' = '  -0.26%  '
    'D20' = '12.32'
    'E20' = '  -1.16%  '
    'E21' = '  -0.07%  '
    'D22' = '7.424'
    'E22' = '  -0.36%  '
    'D23' = '1.000'
    'E23' = '  -0.17%  '
    'D24' = '159.03'
    'E24' = '  +1.34%  '
    'D25' = '0.1393'
    'E25' = '  +0.59%  '
    'D26' = '8.429'
    'E26' = '  +0.25%  '
    'E27' = '  -0.33%  '
    'D28' = '1.431'
    'E28' = '  +8.19%  '
    'D29' = '1.473'
    'E29' = '  +0.11%  '
    'D30' = '0.05657'
    'E30' = '  -0.43%  '
    'D31' = '4.110'
    'E31' = '  -0.60%  '
    'D32' = '4.041'
    'E32' = '  -0.08%  '
    'B33' = 'ARBITRUM'
    'C33' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D33' = '1.153'
    'E33' = '  -0.93%  '
    'B34' = 'LidoDAOToken'
    'C34' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D34' = '1.820'
    'E34' = '  -1.47%  '
    'D35' = '0.6966'
    'E35' = '  -1.73%  '
    'D36' = '2.586'
    'E36' = '  +0.01%  '
    'E37' = '  +1.79%  '
    'D38' = '1.237.07'
    'E38' = '  +1.31%  '
    'E39' = '  -2.27%  '
    'D40' = '6.386'
    'E40' = '  -2.65%  '
    'D41' = '0.9012'
    'E41' = '  -1.01%  '
    'D42' = '0.9996'
    'E42' = '  -0.20%  '
    'B43' = 'RocketPoolETH'
    'C43' = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
    'D43' = '2.005.77'
    'E43' = '  -0.64%  '
    'B44' = 'Quant'
    'C44' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D44' = '101.30'
    'E44' = '  -0.51%  '
    'B45' = 'Aave'
    'C45' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D45' = '65.41'
    'E45' = '  -1.46%  '
    'B46' = 'Aptos'
    'C46' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D46' = '7.133'
    'E46' = '  +0.13%  '
    'B47' = 'BabyDogeCoin'
    'C47' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D47' = '0.00000000117'
    'E47' = '  -4.04%  '
    'B48' = 'TheSandbox'
    'C48' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D48' = '0.3996'
    'E48' = '  -0.69%  '
    'B49' = 'Algorand'
    'C49' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D49' = '0.1153'
    'E49' = '  +1.02%  '
    'B50' = 'EnergySwap'
    'C50' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D50' = '8.977'
    'E50' = '  -0.40%  '
    'B51' = 'RenderToken'
    'C51' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D51' = '1.671'
    'E51' = '  -0.78%  '
}

foreach ($addr in $cellUpdates.Keys) {
    $col = $addr -replace '[0-9]+$', ''
    $rng = $ws.Range($addr)
    $value = $cellUpdates[$addr]

    if ($col -eq 'D' -or $col -eq 'E') {
        # Force text storage so numeric-looking strings (prices,
        # percentages) are not auto-converted into real numbers.
        $rng.NumberFormat = '@'
        $rng.Value = $value
        $rng.Style = 'Normal'
    } else {
        # Plain text columns (coin name / link) are safe to set directly.
        $rng.Value = $value
    }
}
